$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting rows 10-12 down to 11-13
$ws.Rows.Item(10).Insert()

# New row 10: id=9, adm0="United States", adm1="Philadelphia"
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "United States"
$ws.Cells.Item(10, 3).Value = "Philadelphia"

# Update the id column for the shifted rows (11, 12, 13) since previously
# they held values 9, 10, 11 and now should hold 10, 11, 12
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(13, 1).Value = 12

# Adjust column C width (now not best-fit, custom width)
$ws.Columns.Item(3).ColumnWidth = 11.0

# Update selection to E7
$ws.Range("E7").Select()
